# Update TPM-derived NATMI ligand-receptor metrics (Fgf22-Fgfrl1) with
# recomputed values following the script update referenced in the commit
# message ("update scripts wuth new tpm").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("I2").Value = 0.3768554877997065
$ws.Range("J2").Value = 0.3768554877997065
$ws.Range("M2").Value = 1.193633666666667
$ws.Range("N2").Value = 3.580901
$ws.Range("O2").Value = 0.0852504197348203
$ws.Range("P2").Value = 0.08525041973482028
$ws.Range("Q2").Value = 0.176792663271
$ws.Range("R2").Value = 1.591133969439
$ws.Range("S2").Value = 0.03212708851429543
$ws.Range("T2").Value = 0.03212708851429542

# Row 3
$ws.Range("I3").Value = 0.3768554877997065
$ws.Range("J3").Value = 0.3768554877997065
$ws.Range("O3").Value = 0.6175422122064692
$ws.Range("P3").Value = 0.6175422122064691
$ws.Range("S3").Value = 0.2327241716179788
$ws.Range("T3").Value = 0.2327241716179788

# Row 4
$ws.Range("I4").Value = 0.3768554877997065
$ws.Range("J4").Value = 0.3768554877997065
$ws.Range("M4").Value = 4.028899666666667
$ws.Range("N4").Value = 12.086699
$ws.Range("O4").Value = 0.2877477380576656
$ws.Range("P4").Value = 0.2877477380576655
$ws.Range("Q4").Value = 0.596732416329
$ws.Range("R4").Value = 5.370591746961
$ws.Range("S4").Value = 0.1084393141889837
$ws.Range("T4").Value = 0.1084393141889837

# Row 5
$ws.Range("I5").Value = 0.3768554877997065
$ws.Range("J5").Value = 0.3768554877997065
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.132449
$ws.Range("N5").Value = 0.397347
$ws.Range("O5").Value = 0.009459630001044888
$ws.Range("P5").Value = 0.009459630001044887
$ws.Range("Q5").Value = 0.019617418737
$ws.Range("R5").Value = 0.176556768633
$ws.Range("S5").Value = 0.003564913478448509
$ws.Range("T5").Value = 0.003564913478448509

# Row 6
$ws.Range("G6").Value = 0.2449103333333333
$ws.Range("H6").Value = 0.734731
$ws.Range("I6").Value = 0.6231445122002934
$ws.Range("J6").Value = 0.6231445122002934
$ws.Range("M6").Value = 1.193633666666667
$ws.Range("N6").Value = 3.580901
$ws.Range("O6").Value = 0.0852504197348203
$ws.Range("P6").Value = 0.08525041973482028
$ws.Range("Q6").Value = 0.2923332191812222
$ws.Range("R6").Value = 2.630998972631
$ws.Range("S6").Value = 0.05312333122052487
$ws.Range("T6").Value = 0.05312333122052485

# Row 7
$ws.Range("G7").Value = 0.2449103333333333
$ws.Range("H7").Value = 0.734731
$ws.Range("I7").Value = 0.6231445122002934
$ws.Range("J7").Value = 0.6231445122002934
$ws.Range("O7").Value = 0.6175422122064692
$ws.Range("P7").Value = 0.6175422122064691
$ws.Range("Q7").Value = 2.117621279005556
$ws.Range("R7").Value = 19.05859151105
$ws.Range("S7").Value = 0.3848180405884903
$ws.Range("T7").Value = 0.3848180405884903

# Row 8
$ws.Range("G8").Value = 0.2449103333333333
$ws.Range("H8").Value = 0.734731
$ws.Range("I8").Value = 0.6231445122002934
$ws.Range("J8").Value = 0.6231445122002934
$ws.Range("M8").Value = 4.028899666666667
$ws.Range("N8").Value = 12.086699
$ws.Range("O8").Value = 0.2877477380576656
$ws.Range("P8").Value = 0.2877477380576655
$ws.Range("Q8").Value = 0.9867191603298889
$ws.Range("R8").Value = 8.880472442969
$ws.Range("S8").Value = 0.1793084238686818
$ws.Range("T8").Value = 0.1793084238686818

# Row 9
$ws.Range("G9").Value = 0.2449103333333333
$ws.Range("H9").Value = 0.734731
$ws.Range("I9").Value = 0.6231445122002934
$ws.Range("J9").Value = 0.6231445122002934
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.132449
$ws.Range("N9").Value = 0.397347
$ws.Range("O9").Value = 0.009459630001044888
$ws.Range("P9").Value = 0.009459630001044887
$ws.Range("Q9").Value = 0.03243812873966667
$ws.Range("R9").Value = 0.291943158657
$ws.Range("S9").Value = 0.005894716522596378
$ws.Range("T9").Value = 0.005894716522596377
